# Add the 2022-Q3 quarterly sheet and update the "总计" (totals) summary
# sheet accordingly.
#
# Target layout after the edit:
#   总计 (totals) | 2022-Q3 (NEW) | 2022-Q2 | 2022-Q1 | 2021-Q4 | 2021-Q3
#
# The existing quarter sheets (2022-Q2 / 2022-Q1 / 2021-Q4 / 2021-Q3) keep
# their own data untouched -- they just shift one tab to the right to make
# room for the new "2022-Q3" sheet, which is inserted right after "总计".

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force a (numeric-looking) string to be stored as text rather than
    # letting Excel auto-coerce it to a number, then strip the temporary
    # "Text" number-format back off so the cell ends up plain/unstyled --
    # matching the look of the surrounding data cells.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q3" sheet.
#    Easiest high-fidelity way to match the existing quarter-sheet look
#    (header styling, borders, column widths, etc.) is to clone the
#    "2022-Q2" sheet -- which shares the exact same template -- place the
#    clone immediately before it, rename it, and then overwrite the cell
#    values and trim the extra fund rows it doesn't need.
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2022-Q2")
$templateSheet.Copy($templateSheet)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The template ("2022-Q2") has 4 fund rows (rows 2-5); "2022-Q3" only needs
# 2 (rows 2-3), so drop the two extra fund rows -- originally rows 3 and 5;
# delete the higher row index first so the second delete still targets the
# intended original row.
$q3.Rows.Item(5).Delete()
$q3.Rows.Item(3).Delete()

# Overwrite the remaining fund rows with the 2022-Q3 figures.
$q3.Range("A2").Value = 0
Set-TextValue $q3.Range("B2") "009686"
$q3.Range("C2").Value = "华夏磐利一年定期开放混合A"
Set-TextValue $q3.Range("D2") "10.76"
Set-TextValue $q3.Range("E2") "64.78"
Set-TextValue $q3.Range("F2") "2.51"
Set-TextValue $q3.Range("G2") "0.2701"
$q3.Range("H2").Value = 7

$q3.Range("A3").Value = 1
Set-TextValue $q3.Range("B3") "009687"
$q3.Range("C3").Value = "华夏磐利一年定期开放混合C"
Set-TextValue $q3.Range("D3") "0.43"
Set-TextValue $q3.Range("E3") "64.78"
Set-TextValue $q3.Range("F3") "2.51"
Set-TextValue $q3.Range("G3") "0.0108"
$q3.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a fresh row for 2022-Q3 right under
#    the header and renumber the existing rows beneath it.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Rows.Item(2).Insert()
# Insert() interpolates a border/format for the new row; strip that back
# down to a plain (unstyled) data row like every other data row here.
$totals.Range("B2:D2").ClearFormats()
# ... except column A, which does carry the header-matching style in every
# data row -- clone it from the row right below.
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 2
$totals.Range("D2").Value = 0.28

$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
$totals.Range("A6").Value = 4

# ---------------------------------------------------------------------
# 3. Restore the originally-active tab ("2021-Q3") -- cloning/renaming
#    sheets above left the new "2022-Q3" sheet selected instead.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
